$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.164.21"
$ws.Range("E2").Value = "  +0.82%  "
$ws.Range("D3").Value = "1.901.59"
$ws.Range("E3").Value = "  +1.36%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").Value = "'306.19"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("D7").Value = "'0.5249"
$ws.Range("E7").Value = "  +1.64%  "
$ws.Range("D8").Value = "'0.3773"
$ws.Range("E8").Value = "  +1.69%  "
$ws.Range("D9").Value = "'0.07256"
$ws.Range("E9").Value = "  +1.02%  "
$ws.Range("D10").Value = "'21.15"
$ws.Range("E10").Value = "  +2.07%  "
$ws.Range("D11").Value = "'0.8995"
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("D12").Value = "'0.08363"
$ws.Range("E12").Value = "  +10.68%  "
$ws.Range("D13").Value = "1.899.12"
$ws.Range("E13").Value = "  +1.17%  "
$ws.Range("D14").Value = "'94.75"
$ws.Range("E14").Value = "  +0.04%  "
$ws.Range("E15").Value = "  +0.39%  "
$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "  +0.26%  "
$ws.Range("D17").Value = "'0.000008592"
$ws.Range("E17").Value = "  +1.31%  "
$ws.Range("D18").Value = "'14.51"
$ws.Range("E18").Value = "  +1.83%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("D20").Value = "27.201.61"
$ws.Range("E20").Value = "  +0.80%  "
$ws.Range("D21").Value = "'5.063"
$ws.Range("E21").Value = "  +0.73%  "
$ws.Range("D22").Value = "2.141.85"
$ws.Range("E22").Value = "  +2.10%  "
$ws.Range("D23").Value = "'10.59"
$ws.Range("E23").Value = "  +1.87%  "
$ws.Range("D24").Value = "'6.424"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").Value = "'2.287"
$ws.Range("E25").Value = "  +8.40%  "
$ws.Range("D26").Value = "'146.63"
$ws.Range("E26").Value = "  +0.54%  "
$ws.Range("D27").Value = "'1.758"
$ws.Range("E27").Value = "  -1.38%  "
$ws.Range("E28").Value = "  +0.66%  "
$ws.Range("D29").Value = "'114.76"
$ws.Range("E29").Value = "  +0.27%  "
$ws.Range("D30").Value = "'4.928"
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("D31").Value = "'4.786"
$ws.Range("E31").Value = "  +0.74%  "
$ws.Range("D32").Value = "'0.09230"
$ws.Range("E32").Value = "  +0.51%  "
$ws.Range("D33").Value = "'0.8146"
$ws.Range("E33").Value = "  +8.16%  "
$ws.Range("D34").Value = "'0.05051"
$ws.Range("E34").Value = "  +0.30%  "
$ws.Range("D35").Value = "'1.237"
$ws.Range("E35").Value = "  +5.40%  "
$ws.Range("D36").Value = "'2.965"
$ws.Range("E36").Value = "  -1.19%  "
$ws.Range("D37").Value = "'3.358"
$ws.Range("E37").Value = "  +3.09%  "
$ws.Range("E38").Value = "  +2.89%  "
$ws.Range("D39").Value = "'0.5689"
$ws.Range("E39").Value = "  +1.91%  "
$ws.Range("D40").Value = "'0.01973"
$ws.Range("E40").Value = "  -0.90%  "
$ws.Range("D41").Value = "'1.075"
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("D42").Value = "'6.661"
$ws.Range("E42").Value = "  +1.37%  "
$ws.Range("D43").Value = "'8.957"
$ws.Range("E43").Value = "  +2.46%  "
$ws.Range("D44").Value = "'118.54"
$ws.Range("E44").Value = "  +1.79%  "
$ws.Range("D45").Value = "'0.1511"
$ws.Range("E45").Value = "  +0.62%  "
$ws.Range("D46").Value = "'0.4825"
$ws.Range("E46").Value = "  +1.05%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'10.20"
$ws.Range("E47").Value = "  +0.59%  "
$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D48").Value = "'1.002"
$ws.Range("E48").Value = "  +0.33%  "
$ws.Range("D49").Value = "'1.609"
$ws.Range("E49").Value = "  +2.77%  "
$ws.Range("D50").Value = "'37.46"
$ws.Range("E50").Value = "  +0.93%  "
$ws.Range("D51").Value = "'63.52"
$ws.Range("E51").Value = "  +0.29%  "
